# Insert a new row at position 943, shifting all existing data rows
# (old rows 943-1013) down by one (to 944-1014), matching the weekly
# price update for "Fruta / hortaliza, semanal".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(943).Insert()

$ws.Range("A943").Value = 4
$ws.Range("B943").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C943").Value = "Los Lagos"
$ws.Range("D943").Value = 45265
$ws.Range("E943").Value = 10
$ws.Range("F943").Value = "Fruta"
$ws.Range("G943").Value = 100108
$ws.Range("H943").Value = "Tropicales y subtropicales"
$ws.Range("I943").Value = 100108006
$ws.Range("J943").Value = "Plátano"
$ws.Range("K943").Value = "Sin especificar"
$ws.Range("L943").Value = "Primera Pintón"
$ws.Range("M943").Value = 900
$ws.Range("N943").Value = 29000
$ws.Range("O943").Value = 30000
$ws.Range("P943").Value = 29556
$ws.Range("Q943").Value = "$/caja 20 kilos"
$ws.Range("R943").Value = "Ecuador"
$ws.Range("S943").Value = 1478
$ws.Range("T943").Value = 20
